$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.044.92"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.684.76"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.58%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.920.36"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "1.671.37"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "27.070.84"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "0.0₃0736"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.54%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D33").Value = "1.511.71"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("E35").Value = "  +4.43%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.591"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.922"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "1.823.56"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0507"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
